$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove legacy ArcGIS columns: OBJECTID_1, OBJECTID__, Shape__ (old A:C)
$ws.Columns("A:C").Delete()

# After the above shift, the remaining legacy trailer columns
# (Shape_Length, Shape_Area, OBJECTID) now sit at I:K - remove them too,
# leaving Area_SqKm as the last column.
$ws.Columns("I:K").Delete()

# Clear the literal "<Null>" placeholders left behind in empty fields
$ws.Range("C2").ClearContents()
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("H2").ClearContents()

# Freeze the header row
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# Set readable column widths for the remaining fields
$ws.Columns("A").ColumnWidth = 12.877604166666666
$ws.Columns("B").ColumnWidth = 10.877604166666666
$ws.Columns("C").ColumnWidth = 12.877604166666666
$ws.Columns("D:E").ColumnWidth = 16.877604166666668
$ws.Columns("F:G").ColumnWidth = 9.877604166666666
$ws.Columns("H").ColumnWidth = 6.877604166666667
$ws.Columns("I").ColumnWidth = 8.877604166666666
